# Weekly update: a new price record (week of 2022-03-22, serial 44642) is
# inserted as a new data row right before the existing row that is currently
# on row 406. Excel's native "insert row" shifts every row from 406 downward
# by one (so the old row 512 becomes row 513, and the sheet's used range
# grows from A1:R512 to A1:R513); we then populate the freshly inserted row
# with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 406 - this pushes the former rows
# 406..512 down to 407..513 and extends the sheet dimension accordingly.
$ws.Rows.Item(406).Insert()

# Populate the new row 406 with the new weekly record.
$ws.Range("A406").Value = 6
$ws.Range("B406").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C406").Value = "Metropolitana"
$ws.Range("D406").Value = 44642
$ws.Range("E406").Value = 13
$ws.Range("F406").Value = 100112012
$ws.Range("G406").Value = "Espinaca"
$ws.Range("H406").Value = "Sin especificar"
$ws.Range("I406").Value = "Primera"
$ws.Range("J406").Value = 590
$ws.Range("K406").Value = 6000
$ws.Range("L406").Value = 6500
$ws.Range("M406").Value = 6220
$ws.Range("N406").Value = "`$/cuna 10 kilos"
$ws.Range("O406").Value = "Región Metropolitana"
$ws.Range("P406").Value = 622
$ws.Range("Q406").Value = 10
$ws.Range("R406").Value = "Hortaliza"
